$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1154.6666
$ws.Range("I12").Value = 785.6
$ws.Range("J12").Value = 3000
$ws.Range("K12").Value = 785.6
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -615.6
$ws.Range("N12").Value = -3340

$ws.Range("H40").Value = 7043.3335
$ws.Range("I40").Value = 3800
$ws.Range("J40").Value = 7854.1665
$ws.Range("K40").Value = 3800
$ws.Range("L40").Value = 7854.1665
$ws.Range("M40").Value = -3625
$ws.Range("N40").Value = -8204.166499999999

$ws.Range("H70").Value = 3026.2856
$ws.Range("I70").Value = 2724.5
$ws.Range("K70").Value = 8173.5
$ws.Range("M70").Value = -7903.5

$ws.Range("H73").Value = 3026.2856
$ws.Range("I73").Value = 2724.5
$ws.Range("K73").Value = 8173.5
$ws.Range("M73").Value = -7237.5

$ws.Range("H87").Value = 67372.375
$ws.Range("J87").Value = 67372.375
$ws.Range("L87").Value = 67372.375
$ws.Range("N87").Value = -69868.375

$ws.Range("H90").Value = 67372.375
$ws.Range("J90").Value = 67372.375
$ws.Range("L90").Value = 202117.125
$ws.Range("N90").Value = -214597.125

$ws.Range("H100").Value = 2215.2
$ws.Range("J100").Value = 2596.3333
$ws.Range("L100").Value = 2596.3333
$ws.Range("N100").Value = -3678.3333

$ws.Range("H133").Value = 100039.5
$ws.Range("J133").Value = 101184.336
$ws.Range("L133").Value = 101184.336
$ws.Range("N133").Value = -111304.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 35714936
$ws.Range("I110").Value = 45455160
$ws.Range("J110").Value = 795.6667
$ws.Range("K110").Value = 45455160
$ws.Range("L110").Value = 795.6667
$ws.Range("M110").Value = -45453115
$ws.Range("N110").Value = -4885.6667

$ws.Range("H132").Value = 2171.325
$ws.Range("I132").Value = 2339.5667
$ws.Range("J132").Value = 1666.6
$ws.Range("K132").Value = 7018.7001
$ws.Range("L132").Value = 4999.799999999999
$ws.Range("M132").Value = -4488.7001
$ws.Range("N132").Value = -10059.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 993.4286
$ws.Range("I22").Value = 990.2
$ws.Range("J22").Value = 1001.5
$ws.Range("K22").Value = 990.2
$ws.Range("L22").Value = 1001.5
$ws.Range("M22").Value = -817.2
$ws.Range("N22").Value = -1347.5

$ws.Range("H81").Value = 22941.908
$ws.Range("J81").Value = 22941.908
$ws.Range("L81").Value = 22941.908
$ws.Range("N81").Value = -25063.908

$ws.Range("H84").Value = 22941.908
$ws.Range("J84").Value = 22941.908
$ws.Range("L84").Value = 68825.724
$ws.Range("N84").Value = -79433.724

$ws.Range("H99").Value = 2223
$ws.Range("I99").Value = 2122.4443
$ws.Range("J99").Value = 2524.6667
$ws.Range("K99").Value = 2122.4443
$ws.Range("L99").Value = 2524.6667
$ws.Range("M99").Value = -624.4443000000001
$ws.Range("N99").Value = -5520.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 32057608
$ws.Range("I58").Value = 27779980
$ws.Range("J58").Value = 35724148
$ws.Range("K58").Value = 27779980
$ws.Range("L58").Value = 35724148
$ws.Range("M58").Value = -27779777
$ws.Range("N58").Value = -35724554

$ws.Range("H68").Value = 57998.223
$ws.Range("J68").Value = 58998
$ws.Range("L68").Value = 58998
$ws.Range("N68").Value = -60496

$ws.Range("H71").Value = 57998.223
$ws.Range("J71").Value = 58998
$ws.Range("L71").Value = 176994
$ws.Range("N71").Value = -184482

$ws.Range("H74").Value = 75499.5
$ws.Range("J74").Value = 83428
$ws.Range("L74").Value = 83428
$ws.Range("N74").Value = -85176

$ws.Range("H77").Value = 75499.5
$ws.Range("J77").Value = 83428
$ws.Range("L77").Value = 250284
$ws.Range("N77").Value = -259020

$ws.Range("H105").Value = 2214.125
$ws.Range("I105").Value = 2260
$ws.Range("J105").Value = 1893
$ws.Range("K105").Value = 2260
$ws.Range("L105").Value = 1893
$ws.Range("M105").Value = -513
$ws.Range("N105").Value = -5387

$ws.Range("H132").Value = 3189.6086
$ws.Range("I132").Value = 3198.2273
$ws.Range("K132").Value = 9594.6819
$ws.Range("M132").Value = -7064.6819

$ws.Range("H136").Value = 32057608
$ws.Range("I136").Value = 27779980
$ws.Range("J136").Value = 35724148
$ws.Range("K136").Value = 83339940
$ws.Range("L136").Value = 107172444
$ws.Range("M136").Value = -83337390
$ws.Range("N136").Value = -107177544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50853212
$ws.Range("I4").Value = 127126370
$ws.Range("K4").Value = 381379110
$ws.Range("M4").Value = -381378998

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 9750.5
$ws.Range("J7").Value = 9667.333000000001
$ws.Range("L7").Value = 9667.333000000001
$ws.Range("N7").Value = -9891.333000000001

$ws.Range("H8").Value = 9750.5
$ws.Range("J8").Value = 9667.333000000001
$ws.Range("L8").Value = 9667.333000000001
$ws.Range("N8").Value = -9945.333000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 16500
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("N3").Value = -12224

$ws.Range("H15").Value = 16500
$ws.Range("J15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12340

$ws.Range("H22").Value = 5365.364
$ws.Range("J22").Value = 5365.364
$ws.Range("L22").Value = 5365.364
$ws.Range("N22").Value = -5955.364

$ws.Range("H27").Value = 5365.364
$ws.Range("J27").Value = 5365.364
$ws.Range("L27").Value = 5365.364
$ws.Range("N27").Value = -5579.364

$ws.Range("H93").Value = 1847.5385
$ws.Range("I93").Value = 2083
$ws.Range("J93").Value = 1208.4286
$ws.Range("K93").Value = 2083
$ws.Range("L93").Value = 1208.4286
$ws.Range("M93").Value = -835
$ws.Range("N93").Value = -3704.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 400000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H81").Value = 1825.5714
$ws.Range("I81").Value = 1825.5714
$ws.Range("K81").Value = 3651.1428
$ws.Range("M81").Value = -2590.1428

$ws.Range("H84").Value = 1825.5714
$ws.Range("I84").Value = 1825.5714
$ws.Range("K84").Value = 18255.714
$ws.Range("M84").Value = -12951.714
